# Insert a new bulleted "journal" entry right after the entry that ends
# with "...was missed during the first version of language porting)."
# and before the "Version updated to 1.0.14" entry.

$d = $word.ActiveDocument

# Locate the run of text that ends the preceding list item. Find.Execute
# collapses/replaces the range it's called on to the matched text, so we
# search within a fresh Range anchored on the whole story.
$rng = $d.Content
$found = $rng.Find.Execute(
    "was missed during the first version of language porting).",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor text for the new journal entry."
}

# Get the paragraph that contains the match, then collapse a range to its
# very end (just before the paragraph mark) so a new paragraph can be
# inserted immediately after it, inheriting the same list/style formatting.
$anchorPara = $rng.Paragraphs(1)
$endRange = $anchorPara.Range
$endRange.Collapse(0)   # wdCollapseEnd
$endRange.InsertParagraphAfter()

# The newly created (empty) paragraph is now the paragraph right after the
# anchor paragraph; fill in its text.
$newPara = $anchorPara.Next()
$newPara.Range.Text = "Updated sw.js with lang.js (Removed lang.en.js and lang.de.js)"
